# Test Suite1.xlsx edit script
# - "Test Steps" sheet: insert a new row (71) duplicating the hotel
#   "outside-selected-date" step, and repurpose the old row 70 into the
#   "outside field" variant.
# - "Test Data" sheet: add SKIP / PASS / error-message markers in columns
#   E/F, widen those columns to fit, and make "Test Data" the active sheet
#   (selection moves to D6).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Test Steps" sheet
# ---------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Steps")

# Insert a new row at 71 (existing rows 71-76 shift down to 72-77),
# copying the format of row 70 so the new row matches its neighbours.
$wsSteps.Rows.Item(71).Insert()
$wsSteps.Range("A70:E70").Copy()
$wsSteps.Range("A71:E71").PasteSpecial(-4122)

# New row 71: same as the (old) row 70 used to be - the "choose date
# outside the selected range" step.
$wsSteps.Range("A71").Value = "hotelSearchlistTestCase"
$wsSteps.Range("B71").Value = ""
$wsSteps.Range("C71").Value = "click"
$wsSteps.Range("D71").Value = "homepage.hotel.getchlOutselecteddate"
$wsSteps.Range("E71").Value = ""

# Row 70 becomes the new "outside field" step.
$wsSteps.Range("D70").Value = "homepage.hotel.outsidefield"

# ---------------------------------------------------------------------
# 2. "Test Data" sheet
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Test Data")

# Column widths for the new, wider E/F contents (values chosen so the
# saved width lands as close as possible to the target 31.285.../6.57...
# given this engine's column-width quantization).
$wsData.Columns.Item(5).ColumnWidth = 30.45182295
$wsData.Columns.Item(6).ColumnWidth = 5.7369792

# Mark each login-data row with a SKIP flag in column E (rows 3-8).
$wsData.Range("E3").Value = "SKIP"
$wsData.Range("E4").Value = "SKIP"
$wsData.Range("E5").Value = "SKIP"
$wsData.Range("E6").Value = "SKIP"
$wsData.Range("E7").Value = "SKIP"
$wsData.Range("E8").Value = "SKIP"

# Row 12 (hotel search data): expected-error text + PASS flag.
$wsData.Range("E12").Value = "Please enter a valid Email Address"
$wsData.Range("F12").Value = "PASS"

# Selection on "Test Data" moves to D6, and the sheet becomes the active
# (selected) tab of the workbook.
[void]$wsData.Range("D6").Select()
[void]$wsData.Activate()

Write-Host "edit complete"
